# Update team-specific time-matrix probabilities (added team specific time data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1972477064220184
$ws.Range("C2").Value = 0.555045871559633
$ws.Range("J2").Value = 0.01834862385321101
$ws.Range("P2").Value = 0.1345565749235474
$ws.Range("S2").Value = 0.09480122324159021

# Row 3
$ws.Range("B3").Value = 0.007936507936507936
$ws.Range("C3").Value = 0.02380952380952381
$ws.Range("J3").Value = 0.02645502645502645
$ws.Range("P3").Value = 0.7433862433862434
$ws.Range("S3").Value = 0.1984126984126984

# Row 4
$ws.Range("J4").Value = 0.02564102564102564
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.282051282051282

# Row 5
$ws.Range("J5").Value = 0.25
$ws.Range("P5").Value = 0.75

# Row 6
$ws.Range("B6").Value = 0.06954436450839328
$ws.Range("D6").Value = 0.007194244604316547
$ws.Range("F6").Value = 0.04796163069544365
$ws.Range("J6").Value = 0.2661870503597122
$ws.Range("O6").Value = 0.01918465227817746
$ws.Range("Q6").Value = 0.1654676258992806
$ws.Range("R6").Value = 0.07673860911270983
$ws.Range("S6").Value = 0.3477218225419664

# Row 7
$ws.Range("B7").Value = 0.09876543209876543
$ws.Range("D7").Value = 0.01975308641975309
$ws.Range("E7").Value = 0.002469135802469136
$ws.Range("F7").Value = 0.04197530864197531
$ws.Range("J7").Value = 0.1777777777777778
$ws.Range("O7").Value = 0.01728395061728395
$ws.Range("Q7").Value = 0.1703703703703704
$ws.Range("R7").Value = 0.0691358024691358
$ws.Range("S7").Value = 0.4024691358024691

# Row 8
$ws.Range("B8").Value = 0.1192660550458716
$ws.Range("D8").Value = 0.0194954128440367
$ws.Range("F8").Value = 0.09059633027522936
$ws.Range("J8").Value = 0.1077981651376147
$ws.Range("O8").Value = 0.01490825688073395
$ws.Range("Q8").Value = 0.1777522935779816
$ws.Range("R8").Value = 0.08944954128440367
$ws.Range("S8").Value = 0.3807339449541284

# Row 9
$ws.Range("B9").Value = 0.08860759493670886
$ws.Range("D9").Value = 0.006329113924050633
$ws.Range("E9").Value = 0.003164556962025316
$ws.Range("F9").Value = 0.05696202531645569
$ws.Range("J9").Value = 0.129746835443038
$ws.Range("O9").Value = 0.0189873417721519
$ws.Range("Q9").Value = 0.1930379746835443
$ws.Range("R9").Value = 0.06645569620253164
$ws.Range("S9").Value = 0.4367088607594937

# Row 10
$ws.Range("B10").Value = 0.1174721189591078
$ws.Range("D10").Value = 0.01895910780669145
$ws.Range("E10").Value = 0.001486988847583643
$ws.Range("F10").Value = 0.06356877323420074
$ws.Range("J10").Value = 0.128996282527881
$ws.Range("O10").Value = 0.01561338289962825
$ws.Range("Q10").Value = 0.2245353159851301
$ws.Range("R10").Value = 0.08178438661710037
$ws.Range("S10").Value = 0.3475836431226766

# Row 11
$ws.Range("G11").Value = 0.1529968454258675
$ws.Range("J11").Value = 0.08517350157728706
$ws.Range("K11").Value = 0.1971608832807571
$ws.Range("L11").Value = 0.5583596214511041
$ws.Range("S11").Value = 0.006309148264984227

# Row 12
$ws.Range("G12").Value = 0.7131147540983607
$ws.Range("J12").Value = 0.2377049180327869
$ws.Range("K12").Value = 0.00273224043715847
$ws.Range("L12").Value = 0.02459016393442623
$ws.Range("S12").Value = 0.02185792349726776

# Row 13
$ws.Range("F13").Value = 0.01136363636363636
$ws.Range("G13").Value = 0.6363636363636364
$ws.Range("J13").Value = 0.2840909090909091
$ws.Range("S13").Value = 0.06818181818181818

# Row 15
$ws.Range("F15").Value = 0.01025641025641026
$ws.Range("H15").Value = 0.1358974358974359
$ws.Range("I15").Value = 0.0641025641025641
$ws.Range("J15").Value = 0.382051282051282
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("M15").Value = 0.02307692307692308
$ws.Range("O15").Value = 0.04871794871794872
$ws.Range("S15").Value = 0.2692307692307692

# Row 16
$ws.Range("F16").Value = 0.01456310679611651
$ws.Range("H16").Value = 0.1601941747572816
$ws.Range("I16").Value = 0.0703883495145631
$ws.Range("J16").Value = 0.4344660194174757
$ws.Range("K16").Value = 0.0970873786407767
$ws.Range("M16").Value = 0.01941747572815534
$ws.Range("N16").Value = 0.002427184466019417
$ws.Range("O16").Value = 0.03398058252427184
$ws.Range("S16").Value = 0.1674757281553398

# Row 17
$ws.Range("F17").Value = 0.01048218029350105
$ws.Range("H17").Value = 0.1865828092243187
$ws.Range("I17").Value = 0.06708595387840671
$ws.Range("J17").Value = 0.4412997903563941
$ws.Range("K17").Value = 0.09748427672955975
$ws.Range("M17").Value = 0.01781970649895178
$ws.Range("O17").Value = 0.05031446540880503
$ws.Range("S17").Value = 0.1289308176100629

# Row 18
$ws.Range("F18").Value = 0.01066666666666667
$ws.Range("H18").Value = 0.1866666666666667
$ws.Range("I18").Value = 0.09066666666666667
$ws.Range("J18").Value = 0.3946666666666667
$ws.Range("K18").Value = 0.104
$ws.Range("M18").Value = 0.02133333333333333
$ws.Range("O18").Value = 0.08266666666666667
$ws.Range("S18").Value = 0.1093333333333333

# Row 19
$ws.Range("F19").Value = 0.01488933601609658
$ws.Range("H19").Value = 0.2056338028169014
$ws.Range("I19").Value = 0.06680080482897384
$ws.Range("J19").Value = 0.3907444668008048
$ws.Range("K19").Value = 0.1223340040241449
$ws.Range("M19").Value = 0.02092555331991952
$ws.Range("N19").Value = 0.0004024144869215292
$ws.Range("O19").Value = 0.05955734406438632
$ws.Range("S19").Value = 0.1187122736418511
